$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1450
$ws.Range("I28").Value = 1450
$ws.Range("K28").Value = 1450
$ws.Range("M28").Value = -965

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8999.736999999999
$ws.Range("J40").Value = 9221.888999999999
$ws.Range("L40").Value = 9221.888999999999
$ws.Range("N40").Value = -9571.888999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2093.7144
$ws.Range("I137").Value = 2097.8333
$ws.Range("J137").Value = 2069
$ws.Range("K137").Value = 6293.499899999999
$ws.Range("L137").Value = 6207
$ws.Range("M137").Value = -3743.499899999999
$ws.Range("N137").Value = -11307

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3152.4
$ws.Range("I138").Value = 1177
$ws.Range("K138").Value = 3531
$ws.Range("M138").Value = 1609

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3680.9546
$ws.Range("I32").Value = 3680.9546
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3680.9546
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3393.9546
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 7809.5
$ws.Range("I122").Value = 9269
$ws.Range("K122").Value = 27807
$ws.Range("M122").Value = -25357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3598.8
$ws.Range("I132").Value = 3748.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11245.5
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -8715.5
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8563.625
$ws.Range("I134").Value = 5299.4
$ws.Range("J134").Value = 14004
$ws.Range("K134").Value = 15898.2
$ws.Range("L134").Value = 42012
$ws.Range("M134").Value = -13363.2
$ws.Range("N134").Value = -47082

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6919.846
$ws.Range("I31").Value = 8209.286
$ws.Range("J31").Value = 5415.5
$ws.Range("K31").Value = 8209.286
$ws.Range("L31").Value = 5415.5
$ws.Range("M31").Value = -7914.286
$ws.Range("N31").Value = -6005.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6919.846
$ws.Range("I34").Value = 8209.286
$ws.Range("J34").Value = 5415.5
$ws.Range("K34").Value = 8209.286
$ws.Range("L34").Value = 5415.5
$ws.Range("M34").Value = -8007.286
$ws.Range("N34").Value = -5819.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2486.3333
$ws.Range("I58").Value = 2486.3333
$ws.Range("K58").Value = 2486.3333
$ws.Range("M58").Value = -2283.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 20000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 20000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -19489
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10425
$ws.Range("J62").Value = 7897.5
$ws.Range("L62").Value = 7897.5
$ws.Range("N62").Value = -9145.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 10425
$ws.Range("J65").Value = 7897.5
$ws.Range("L65").Value = 39487.5
$ws.Range("N65").Value = -45727.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 989.8
$ws.Range("I99").Value = 937.25
$ws.Range("K99").Value = 937.25
$ws.Range("M99").Value = 560.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1496.5
$ws.Range("I105").Value = 1496.5
$ws.Range("K105").Value = 1496.5
$ws.Range("M105").Value = 250.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 989.8
$ws.Range("I126").Value = 937.25
$ws.Range("K126").Value = 2811.75
$ws.Range("M126").Value = -341.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2652.75
$ws.Range("I132").Value = 2652.75
$ws.Range("K132").Value = 7958.25
$ws.Range("M132").Value = -5428.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2486.3333
$ws.Range("I136").Value = 2486.3333
$ws.Range("K136").Value = 7458.999899999999
$ws.Range("M136").Value = -4908.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 788.6
$ws.Range("I14").Value = 788.6
$ws.Range("K14").Value = 2365.8
$ws.Range("M14").Value = -2192.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4345.4546
$ws.Range("I131").Value = 2400
$ws.Range("J131").Value = 4777.778
$ws.Range("K131").Value = 7200
$ws.Range("L131").Value = 14333.334
$ws.Range("M131").Value = -2160
$ws.Range("N131").Value = -24413.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5008691
$ws.Range("I36").Value = 6668254.5
$ws.Range("K36").Value = 6668254.5
$ws.Range("M36").Value = -6667769.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 400204.5
$ws.Range("I58").Value = 400204.5
$ws.Range("K58").Value = 400204.5
$ws.Range("M58").Value = -399927.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7070.2856
$ws.Range("I80").Value = 6698.6
$ws.Range("K80").Value = 6698.6
$ws.Range("M80").Value = -5700.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 7070.2856
$ws.Range("I83").Value = 6698.6
$ws.Range("K83").Value = 33493
$ws.Range("M83").Value = -28501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6038
$ws.Range("I122").Value = 4829.2856
$ws.Range("J122").Value = 14499
$ws.Range("K122").Value = 14487.8568
$ws.Range("L122").Value = 43497
$ws.Range("M122").Value = -12037.8568
$ws.Range("N122").Value = -48397

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 15000
$ws.Range("I50").Value = 25000
$ws.Range("J50").Value = 5000
$ws.Range("K50").Value = 25000
$ws.Range("L50").Value = 5000
$ws.Range("M50").Value = -24363
$ws.Range("N50").Value = -6274

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 60000000
$ws.Range("I56").Value = 60000000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 60000000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -59999309
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3884.3333
$ws.Range("I132").Value = 3884.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11652.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9122.999899999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 55000
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 20000
$ws.Range("I40").Value = 20000
$ws.Range("K40").Value = 20000
$ws.Range("M40").Value = -19851

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 9070
$ws.Range("I51").Value = 9070
$ws.Range("K51").Value = 9070
$ws.Range("M51").Value = -8560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4197.4614
$ws.Range("I132").Value = 1730
$ws.Range("K132").Value = 5190
$ws.Range("M132").Value = -2660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1512.5
$ws.Range("I136").Value = 1512.5
$ws.Range("K136").Value = 4537.5
$ws.Range("M136").Value = -1987.5
